# NuvoControl_0011_StundenNachweis_Grafisch.xlsx
# Commit: "Vorbereitung Treffen mit Hans Rudin - Stundenachweis versioniert und
#          PDF erstellt (Rev. C) - Stundennachweis grafisch erfasst und PDF
#          erstellt (Rev. C) - Dokumenten Status Liste nachgeführt."
#
# Concretely: fill in the weekly hours (rows 16-19) on "Eingabedaten" for the
# newly logged weeks, and append a new "Revision C" line on the "Revision"
# sheet. Every other part of the diff (chart caches, the computed sheet
# "Eingabedaten (berechnet)", the various SUM() totals, chart axis ids /
# page-margins ...) is purely a downstream consequence of these inputs and of
# Excel re-saving the workbook, so it is left to the engine's own recalc.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Eingabedaten" - enter the hours booked during weeks 11-14 (rows 16-19)
# ---------------------------------------------------------------------------
$dataSheet = $wb.Worksheets.Item("Eingabedaten")

# Week 39930-39936 (row 16)
$dataSheet.Range("E16").Formula = "=0.5"
$dataSheet.Range("F16").Formula = "=1.5+2+2.5+1"
$dataSheet.Range("G16").Formula = "=0.5"
$dataSheet.Range("I16").Formula = "=2+0.25+1"
$dataSheet.Range("L16").Formula = "=1.5+2.5+1"
$dataSheet.Range("O16").Formula = "=2"

# Week 39937-39943 (row 17)
$dataSheet.Range("F17").Formula = "=3+1.5"
$dataSheet.Range("G17").Formula = "=1+4.5+2"
$dataSheet.Range("H17").Formula = "=1"
$dataSheet.Range("I17").Formula = "=1.75+4"
$dataSheet.Range("L17").Formula = "=3"
$dataSheet.Range("O17").Formula = "=1.75"

# Week 39944-39950 (row 18)
$dataSheet.Range("F18").Formula = "=1+2.5+3.5"
$dataSheet.Range("G18").Formula = "=3+0.5+1.5+0.5+1.5"
$dataSheet.Range("I18").Formula = "=0.5+0.5"
$dataSheet.Range("L18").Formula = "=1+2+4"
$dataSheet.Range("M18").Formula = "=3+2+2"

# Week 39951-39957 (row 19)
$dataSheet.Range("G19").Formula = "=3+3+5+1.5+6+7+2.5"
$dataSheet.Range("I19").Formula = "=1"
$dataSheet.Range("L19").Formula = "=2.5"
$dataSheet.Range("M19").Formula = "=2+3+4+3"
$dataSheet.Range("N19").Formula = "=1"

# Leave the on-sheet selection where the last bit of data was typed in.
$dataSheet.Range("I19").Select()

# ---------------------------------------------------------------------------
# 2) "Revision" - log revision "c" (sent to Hans Rudin for the interim review)
# ---------------------------------------------------------------------------
$revisionSheet = $wb.Worksheets.Item("Revision")
$revisionSheet.Range("C6").Value = "c / 24-Mai-2009"
$revisionSheet.Range("E6").Value = "Version an Hans Rudin Zwecks Zwischenbesprechung"

# ---------------------------------------------------------------------------
# 3) Restore "Revision" as the active sheet/selection (it was active before
#    the edit, and stays the active tab afterwards).
# ---------------------------------------------------------------------------
$revisionSheet.Activate()
$revisionSheet.Range("E6").Select()

$excel.CalculateFull()
